$d = $word.ActiveDocument

# Locate the title paragraph that currently reads "Huella de Carbono223"
# (three runs: "Huella de Carbono", "2", "23") and remove just the
# trailing "23" run, restoring "Huella de Carbono2".
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Huella de Carbono223*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph containing 'Huella de Carbono223'"
}

$r = $target.Range
$txt = $r.Text
$relIdx = $txt.IndexOf("23", $txt.IndexOf("Huella de Carbono") + "Huella de Carbono2".Length - 1)

$absStart = $r.Start + $relIdx
$absEnd = $absStart + 2

$delRange = $d.Range($absStart, $absEnd)
if ($delRange.Text -ne "23") {
    throw "Unexpected range text: $($delRange.Text)"
}
$delRange.Delete()
